# "proxt chrome arguments test"
# Populate column B (rows 2-11) of the LoginData sheet with a list of
# proxy host:port strings, then move the active selection to C7.
#
# Cells are written in the same order the shared-string table ends up
# needing them so the resulting xl/sharedStrings.xml unique-string order
# matches the target workbook (B3, B5, B10, B2, B4, B6, B7, B8, B9, B11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

$ws.Cells.Item(3, 2).Value  = "51.159.207.156:3128"
$ws.Cells.Item(5, 2).Value  = "200.105.215.18:33630"
$ws.Cells.Item(10, 2).Value = "20.54.56.26:8080"
$ws.Cells.Item(2, 2).Value  = "145.40.77.207:3128"
$ws.Cells.Item(4, 2).Value  = "182.253.172.20:8080"
$ws.Cells.Item(6, 2).Value  = "45.233.67.230:999"
$ws.Cells.Item(7, 2).Value  = "51.159.162.151:80"
$ws.Cells.Item(8, 2).Value  = "190.121.207.58:999"
$ws.Cells.Item(9, 2).Value  = "181.129.49.214:999"
$ws.Cells.Item(11, 2).Value = "51.159.162.151:80"

$ws.Range("C7").Select()
